$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.113.09'
$ws.Range('E2').Value = '  -1.11%  '
$ws.Range('D3').Value = '3.227.67'
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '''528.31'
$ws.Range('E5').Value = '  +2.94%  '
$ws.Range('D6').Value = '''170.24'
$ws.Range('E6').Value = '  -3.81%  '
$ws.Range('E7').Value = '  +1.00%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '3.223.50'
$ws.Range('E9').Value = '  -1.15%  '
$ws.Range('D10').Value = '''0.603'
$ws.Range('E10').Value = '  -1.34%  '
$ws.Range('D11').Value = '''53.00'
$ws.Range('E11').Value = '  -6.32%  '
$ws.Range('E12').Value = '  +3.06%  '
$ws.Range('D13').Value = '''0.0000252'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('E14').Value = '  +1.20%  '
$ws.Range('D15').Value = '3.744.76'
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('E16').Value = '  -1.27%  '
$ws.Range('D17').Value = '3.230.16'
$ws.Range('E17').Value = '  -1.70%  '
$ws.Range('D18').Value = '63.015.48'
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('D19').Value = '''17.16'
$ws.Range('E19').Value = '  +0.80%  '
$ws.Range('D20').Value = '''11.03'
$ws.Range('E20').Value = '  +3.16%  '
$ws.Range('D21').Value = '''0.965'
$ws.Range('E21').Value = '  +3.19%  '
$ws.Range('D22').Value = '''365.86'
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').Value = '''3.76'
$ws.Range('E23').Value = '  +4.30%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '''81.01'
$ws.Range('E24').Value = '  +2.24%  '
$ws.Range('B25').Value = 'RenderToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D25').Value = '''11.14'
$ws.Range('E25').Value = '  +2.98%  '
$ws.Range('D26').Value = '''3.99'
$ws.Range('E26').Value = '  +6.57%  '
$ws.Range('D27').Value = '''6.14'
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('D28').Value = '''2.63'
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('D29').Value = '''11.22'
$ws.Range('E29').Value = '  +0.80%  '
$ws.Range('D30').Value = '''8.17'
$ws.Range('E30').Value = '  -0.58%  '
$ws.Range('D31').Value = '''28.41'
$ws.Range('E31').Value = '  +1.04%  '
$ws.Range('D32').Value = '''637.84'
$ws.Range('E32').Value = '  -1.74%  '
$ws.Range('D33').Value = '''6.41'
$ws.Range('E33').Value = '  -2.76%  '
$ws.Range('D34').Value = '''11.17'
$ws.Range('E34').Value = '  +1.81%  '
$ws.Range('E35').Value = '  +3.51%  '
$ws.Range('D36').Value = '''56.84'
$ws.Range('E36').Value = '  -3.11%  '
$ws.Range('D37').Value = '''1.00'
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('D38').Value = '''36.51'
$ws.Range('E38').Value = '  +2.82%  '
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('D41').Value = '0.0₃0713'
$ws.Range('E41').Value = '  +11.33%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '''0.123'
$ws.Range('E42').Value = '  +0.87%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').Value = '''2.55'
$ws.Range('E43').Value = '  +9.73%  '
$ws.Range('D44').Value = '2.868.34'
$ws.Range('E44').Value = '  +0.97%  '
$ws.Range('D45').Value = '''2.95'
$ws.Range('E45').Value = '  +8.00%  '
$ws.Range('D46').Value = '''2.68'
$ws.Range('E46').Value = '  +3.23%  '
$ws.Range('D47').Value = '''0.0391'
$ws.Range('E47').Value = '  +3.47%  '
$ws.Range('D48').Value = '''3.04'
$ws.Range('E48').Value = '  +4.43%  '
$ws.Range('E49').Value = '  -2.42%  '
$ws.Range('E50').Value = '  +2.17%  '
$ws.Range('D51').Value = '''134.03'
$ws.Range('E51').Value = '  +1.52%  '
